$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "62.241.83"
$ws.Range("E2").Value = "  -2.71%  "

Set-TextValue $ws "D3" "3.177.40"
$ws.Range("E3").Value = "  -4.08%  "

$ws.Range("E4").Value = "  +0.12%  "

Set-TextValue $ws "D5" "587.84"
$ws.Range("E5").Value = "  -1.87%  "

Set-TextValue $ws "D6" "134.90"
$ws.Range("E6").Value = "  -5.95%  "

$ws.Range("E7").Value = "  +0.13%  "

Set-TextValue $ws "D8" "3.173.15"
$ws.Range("E8").Value = "  -4.18%  "

$ws.Range("E9").Value = "  -4.31%  "

$ws.Range("E10").Value = "  -5.38%  "

Set-TextValue $ws "D11" "5.24"
$ws.Range("E11").Value = "  -4.88%  "

Set-TextValue $ws "D12" "0.449"
$ws.Range("E12").Value = "  -5.39%  "

$ws.Range("E13").Value = "  -6.42%  "

$ws.Range("E14").Value = "  -4.80%  "

Set-TextValue $ws "D15" "3.697.83"
$ws.Range("E15").Value = "  -4.19%  "

$ws.Range("E16").Value = "  -1.07%  "

Set-TextValue $ws "D17" "3.178.29"
$ws.Range("E17").Value = "  -4.12%  "

Set-TextValue $ws "D18" "62.286.11"
$ws.Range("E18").Value = "  -2.73%  "

Set-TextValue $ws "D19" "6.55"
$ws.Range("E19").Value = "  -5.30%  "

Set-TextValue $ws "D20" "457.46"
$ws.Range("E20").Value = "  -5.23%  "

$ws.Range("E21").Value = "  -2.85%  "

Set-TextValue $ws "D22" "0.701"
$ws.Range("E22").Value = "  -5.32%  "

Set-TextValue $ws "D23" "7.59"
$ws.Range("E23").Value = "  -5.38%  "

Set-TextValue $ws "D24" "13.26"
$ws.Range("E24").Value = "  -2.37%  "

Set-TextValue $ws "D25" "82.41"
$ws.Range("E25").Value = "  -2.77%  "

Set-TextValue $ws "D26" "1.00"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D27" "2.68"
$ws.Range("E27").Value = "  -3.87%  "

$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D28" "1.00"
$ws.Range("E28").Value = "  -0.08%  "

Set-TextValue $ws "D29" "6.89"
$ws.Range("E29").Value = "  -5.38%  "

$ws.Range("E30").Value = "  -4.89%  "

Set-TextValue $ws "D31" "2.03"
$ws.Range("E31").Value = "  -5.91%  "

$ws.Range("E32").Value = "  -7.93%  "

$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("E34").Value = "  -6.59%  "

Set-TextValue $ws "D35" "1.04"
$ws.Range("E35").Value = "  -6.60%  "

Set-TextValue $ws "D36" "5.77"
$ws.Range("E36").Value = "  -3.88%  "

Set-TextValue $ws "D37" "51.04"
$ws.Range("E37").Value = "  -4.26%  "

Set-TextValue $ws "D38" "0.0₃0684"
$ws.Range("E38").Value = "  -9.96%  "

Set-TextValue $ws "D39" "0.0384"
$ws.Range("E39").Value = "  -4.13%  "

Set-TextValue $ws "D40" "2.948.35"
$ws.Range("E40").Value = "  -3.09%  "

Set-TextValue $ws "D41" "407.97"
$ws.Range("E41").Value = "  -5.68%  "

$ws.Range("E42").Value = "  +1.48%  "

Set-TextValue $ws "D43" "2.64"
$ws.Range("E43").Value = "  -4.98%  "

Set-TextValue $ws "D44" "7.99"
$ws.Range("E44").Value = "  -5.39%  "

$ws.Range("E45").Value = "  -7.13%  "

$ws.Range("E46").Value = "  -3.94%  "

Set-TextValue $ws "D48" "35.77"
$ws.Range("E48").Value = "  +0.40%  "

Set-TextValue $ws "D49" "25.43"
$ws.Range("E49").Value = "  -4.37%  "

Set-TextValue $ws "D50" "123.81"
$ws.Range("E50").Value = "  +0.18%  "

$ws.Range("E51").Value = "  -3.90%  "
